# Update cryptos list cell values (prices / 1h volume %) per the
# "Updated cryptos list on Sun Oct 15 09:44:48 UTC 2023 with GitHub Actions" commit.
#
# Columns: A=index(unused), B=Coin, C=Link, D=Price, E=Volume(1h).
# All of D/E (and some B/C, for the two swapped rows) are plain text in the
# source data (not real numbers - many look numeric, e.g. "208.72", "1.01"),
# so force Text number format before assigning to avoid Excel auto-converting
# them into numeric cells / mangling float precision (e.g. "1.62" -> 1.6200000000000001).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "27.058.90"
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = "  +0.55%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.566.45"
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = "  +1.09%  "
$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = "  +0.49%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "208.72"
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = "  +1.20%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.492"
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = "  +0.83%  "
$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = "  +0.43%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "22.12"
$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = "  +0.13%  "
$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = "  +1.27%  "
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = "  +1.95%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.0860"
$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = "  +0.53%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "1.565.05"
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = "  +0.94%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "0.520"
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = "  +0.50%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "27.048.38"
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "61.92"
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.0₃0708"
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = "  +1.59%  "
$ws.Range("B18").Value = "Chainlink"
$ws.Range("C18").Value = "https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link"
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "7.43"
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = "  +2.49%  "
$ws.Range("B19").Value = "BitcoinCash"
$ws.Range("C19").Value = "https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch"
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "216.03"
$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = "  -0.53%  "
$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = "  +0.43%  "
$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = "  +2.42%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "9.20"
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = "  -0.28%  "
$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = "  -0.24%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "154.02"
$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = "  +0.03%  "
$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = "  -0.05%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "15.06"
$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = "  +0.87%  "
$ws.Range("E27").NumberFormat = "@"
$ws.Range("E27").Value = "  +1.83%  "
$ws.Range("E28").NumberFormat = "@"
$ws.Range("E28").Value = "  +0.43%  "
$ws.Range("E29").NumberFormat = "@"
$ws.Range("E29").Value = "  +1.71%  "
$ws.Range("E30").NumberFormat = "@"
$ws.Range("E30").Value = "  +4.37%  "
$ws.Range("E31").NumberFormat = "@"
$ws.Range("E31").Value = "  +0.68%  "
$ws.Range("E32").NumberFormat = "@"
$ws.Range("E32").Value = "  +4.06%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "1.426.34"
$ws.Range("E33").NumberFormat = "@"
$ws.Range("E33").Value = "  +0.27%  "
$ws.Range("E34").NumberFormat = "@"
$ws.Range("E34").Value = "  +12.81%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "1.62"
$ws.Range("E36").NumberFormat = "@"
$ws.Range("E36").Value = "  +2.88%  "
$ws.Range("E37").NumberFormat = "@"
$ws.Range("E37").Value = "  +1.64%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.535"
$ws.Range("E38").NumberFormat = "@"
$ws.Range("E38").Value = "  +1.78%  "
$ws.Range("B39").Value = "FraxShare"
$ws.Range("C39").Value = "https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs"
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "5.82"
$ws.Range("E39").NumberFormat = "@"
$ws.Range("E39").Value = "  +1.99%  "
$ws.Range("B40").Value = "ARBITRUM"
$ws.Range("C40").Value = "https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.813"
$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = "  +0.67%  "
$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = "  +0.47%  "
$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = "  +0.37%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "1.01"
$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = "  +0.66%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "64.93"
$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = "  +0.62%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "1.74"
$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = "  -0.50%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "1.703.08"
$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = "  +1.11%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "86.74"
$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = "  -0.73%  "
$ws.Range("E48").NumberFormat = "@"
$ws.Range("E48").Value = "  +3.37%  "
$ws.Range("E49").NumberFormat = "@"
$ws.Range("E49").Value = "  +0.94%  "
$ws.Range("E50").NumberFormat = "@"
$ws.Range("E50").Value = "  +0.24%  "
$ws.Range("E51").NumberFormat = "@"
$ws.Range("E51").Value = "  +0.48%  "
